$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit turns on iterative calculation for the workbook (calcPr iterate="1").
$excel.Iteration = $true

# The "Methane Blend" market-share table (rows 198-213) lists one row per
# sector for each of two blend technologies ("Biogas" then "Biomass
# Gasification"). The edit adds a new "Iron and Steel" sector row to each
# block, inserted in the same sector order used elsewhere in the table
# (... Petroleum Refining, Iron and Steel, Chemical Products, ...).

# 1) Insert a blank row ahead of "Chemical Products" in the Biogas block.
$ws.Rows.Item(202).Insert()

# 2) Insert a blank row ahead of "Chemical Products" in the Biomass
#    Gasification block (shifted down by one row from the first insert).
$ws.Rows.Item(211).Insert()

# Fill in the new "Iron and Steel" / Biogas row.
$ws.Range("A202").Value = "CIMS.CAN.SK.Iron and Steel.Methane Blend"
$ws.Range("B202").Value = "Service"
$ws.Range("C202").Value = "SK"
$ws.Range("D202").Value = "Iron and Steel"
$ws.Range("E202").Value = "Methane Blend"
$ws.Range("F202").Value = "Biogas"
$ws.Range("G202").Value = "Market share new_max"
$ws.Range("L202").Value = "%"
$ws.Range("M202:W202").Value = 0.01
$ws.Range("X202").Value = "Quick method to simulate ethanol feedstock limits (with FIC) - REPLACE WITH SOMETHING BETTER!"

# Fill in the new "Iron and Steel" / Biomass Gasification row.
$ws.Range("A211").Value = "CIMS.CAN.SK.Iron and Steel.Methane Blend"
$ws.Range("B211").Value = "Service"
$ws.Range("C211").Value = "SK"
$ws.Range("D211").Value = "Iron and Steel"
$ws.Range("E211").Value = "Methane Blend"
$ws.Range("F211").Value = "Biomass Gasification"
$ws.Range("G211").Value = "Market share new_max"
$ws.Range("L211").Value = "%"
$ws.Range("M211:W211").Value = 0.03
$ws.Range("X211").Value = "Quick method to simulate ethanol feedstock limits (with FIC) - REPLACE WITH SOMETHING BETTER!"
